# Update the student Name/Preferences sheet with the refreshed roster:
#  - new student names (Name column, A2:A61)
#  - newly generated, non-duplicated student numbers (Student Number column, B2:B61)
#  - widen column A to fit the longest new name
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is @(Name, StudentNumber). Row order matches the existing
# Stream column (CS for rows 2-37, DS for rows 38-61), which is left untouched.
$students = @(
    @('Keenan Clarke', 943163.0),
    @('Conrad O''Farrell', 215954.0),
    @('Raymond O''Callaghan', 483586.0),
    @('Larry McGrath', 162616.0),
    @('Dillon O''Callaghan', 413535.0),
    @('Benedict O''Connor', 216821.0),
    @('Gavan Sheehan', 242317.0),
    @('Deane Reid', 782668.0),
    @('Eoin O''Rourke', 639215.0),
    @('Albert O''Sullivan', 711912.0),
    @('Eimhin Kennedy', 897165.0),
    @('George O''Reilly', 967415.0),
    @('Kyran White', 864957.0),
    @('Johnny Quinn', 318243.0),
    @('Tiernan Reid', 334652.0),
    @('Rowan Donovan', 411926.0),
    @('William Kenny', 164239.0),
    @('Oran Stewart', 985242.0),
    @('Andre Kennedy', 332325.0),
    @('Feargal Griffin', 618255.0),
    @('Marc O''Donnell', 846263.0),
    @('Tyrone Nolan', 217494.0),
    @('Aran O''Callaghan', 275951.0),
    @('Brain White', 394631.0),
    @('Casey Lynch', 218418.0),
    @('Tyler Brown', 494149.0),
    @('Jason Walsh', 385338.0),
    @('Ultan O''Callaghan', 182343.0),
    @('Conan Regan', 323432.0),
    @('Hugo Fitzgerald', 157144.0),
    @('Eamonn MacNamara', 841774.0),
    @('Leo Byrne', 582614.0),
    @('Zack O''Leary', 112379.0),
    @('Kieran Stewart', 446258.0),
    @('Jonathon Sweeney', 494649.0),
    @('Padhraic Smith', 153199.0),
    @('Caleb O''Keeffe', 431754.0),
    @('Euan O''Leary', 435261.0),
    @('George Burns', 433533.0),
    @('Noel Brady', 242413.0),
    @('Kelan Burns', 314467.0),
    @('Johnathan Byrne', 278833.0),
    @('Emmanuel Maguire', 518575.0),
    @('Kelvin Kane', 335635.0),
    @('Clive Murray', 892753.0),
    @('Damien Murray', 547634.0),
    @('Eamon Moore', 911749.0),
    @('Clayton Daly', 217149.0),
    @('Feargal MacDermott', 139559.0),
    @('Evan O''Callaghan', 729126.0),
    @('Rory Maher', 339635.0),
    @('Harry Foley', 738232.0),
    @('Darren O''Sullivan', 845522.0),
    @('Bernard Hayes', 233691.0),
    @('Louis White', 585966.0),
    @('Greg Flanagan', 419246.0),
    @('Finnan Reid', 142247.0),
    @('Samuel Thompson', 277775.0),
    @('Darragh Stewart', 134182.0),
    @('Mairtin Nolan', 795535.0)
)

$row = 2
foreach ($s in $students) {
    $ws.Cells.Item($row, 1).Value = $s[0]
    $ws.Cells.Item($row, 2).Value = $s[1]
    $row = $row + 1
}

# Column A ("Name") needs to grow to fit the new, longer names (e.g. "Raymond O'Callaghan").
$ws.Columns.Item(1).ColumnWidth = 19.8
